$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range('D2').Value = '27.960.77'
$ws.Range("D2").Style = "Normal"
$ws.Range('E2').Value = '  +3.29%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range('D3').Value = '1.725.52'
$ws.Range("D3").Style = "Normal"
$ws.Range('E3').Value = '  +2.99%  '

$ws.Range('E4').Value = '  -0.21%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range('D5').Value = '218.79'
$ws.Range("D5").Style = "Normal"
$ws.Range('E5').Value = '  +1.66%  '

$ws.Range('E6').Value = '  +1.21%  '

$ws.Range('E7').Value = '  -0.22%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range('D8').Value = '24.03'
$ws.Range("D8").Style = "Normal"
$ws.Range('E8').Value = '  +13.11%  '

$ws.Range('E9').Value = '  +3.82%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range('D10').Value = '0.0634'
$ws.Range("D10").Style = "Normal"
$ws.Range('E10').Value = '  +2.00%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range('D11').Value = '0.0901'
$ws.Range("D11").Style = "Normal"
$ws.Range('E11').Value = '  +2.12%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range('D12').Value = '1.968.96'
$ws.Range("D12").Style = "Normal"
$ws.Range('E12').Value = '  +3.00%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range('D13').Value = '1.721.24'
$ws.Range("D13").Style = "Normal"
$ws.Range('E13').Value = '  +2.69%  '

$ws.Range('E14').Value = '  +3.67%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range('D15').Value = '0.567'
$ws.Range("D15").Style = "Normal"
$ws.Range('E15').Value = '  +5.81%  '

$ws.Range('E16').Value = '  +2.78%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range('D17').Value = '27.908.35'
$ws.Range("D17").Style = "Normal"

$ws.Range("D18").NumberFormat = "@"
$ws.Range('D18').Value = '244.50'
$ws.Range("D18").Style = "Normal"
$ws.Range('E18').Value = '  +2.97%  '

$ws.Range('E19').Value = '  +2.32%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range('D20').Value = '7.88'
$ws.Range("D20").Style = "Normal"
$ws.Range('E20').Value = '  -3.10%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range('D21').Value = '0.999'
$ws.Range("D21").Style = "Normal"
$ws.Range('E21').Value = '  -0.23%  '

$ws.Range('E22').Value = '  +3.96%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range('D23').Value = '9.75'
$ws.Range("D23").Style = "Normal"
$ws.Range('E23').Value = '  +4.75%  '

$ws.Range('E24').Value = '  +0.76%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range('D25').Value = '149.28'
$ws.Range("D25").Style = "Normal"
$ws.Range('E25').Value = '  +1.22%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range('D26').Value = '7.53'
$ws.Range("D26").Style = "Normal"
$ws.Range('E26').Value = '  +4.18%  '

$ws.Range('E27').Value = '  +2.94%  '

$ws.Range('E28').Value = '  +1.93%  '

$ws.Range('E29').Value = '  -0.26%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range('D30').Value = '0.0512'
$ws.Range("D30").Style = "Normal"
$ws.Range('E30').Value = '  +2.84%  '

$ws.Range('E31').Value = '  +1.72%  '

$ws.Range('E32').Value = '  +2.81%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range('D34').Value = '1.492.84'
$ws.Range("D34").Style = "Normal"
$ws.Range('E34').Value = '  -2.42%  '

$ws.Range('E35').Value = '  -1.78%  '

$ws.Range('E36').Value = '  +3.50%  '

$ws.Range('E37').Value = '  +4.88%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range('D38').Value = '2.41'
$ws.Range("D38").Style = "Normal"
$ws.Range('E38').Value = '  +0.31%  '

$ws.Range('E39').Value = '  +0.53%  '

$ws.Range('E40').Value = '  -0.44%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range('D41').Value = '71.54'
$ws.Range("D41").Style = "Normal"
$ws.Range('E41').Value = '  +5.76%  '

$ws.Range('E42').Value = '  +5.72%  '

$ws.Range('E43').Value = '  -0.22%  '

$ws.Range('B44').Value = 'RocketPoolETH'
$ws.Range('C44').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D44").NumberFormat = "@"
$ws.Range('D44').Value = '1.873.55'
$ws.Range("D44").Style = "Normal"
$ws.Range('E44').Value = '  +2.99%  '

$ws.Range('B45').Value = 'MXToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D45").NumberFormat = "@"
$ws.Range('D45').Value = '2.28'
$ws.Range("D45").Style = "Normal"
$ws.Range('E45').Value = '  +0.94%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range('D46').Value = '0.793'
$ws.Range("D46").Style = "Normal"
$ws.Range('E46').Value = '  +1.53%  '

$ws.Range('E47').Value = '  +12.46%  '

$ws.Range('E48').Value = '  +0.52%  '

$ws.Range('E49').Value = '  +3.41%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range('D50').Value = '8.29'
$ws.Range("D50").Style = "Normal"
$ws.Range('E50').Value = '  +3.75%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range('D51').Value = '0.106'
$ws.Range("D51").Style = "Normal"
$ws.Range('E51').Value = '  +1.56%  '
